$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column values (which can look numeric, e.g. "1.003") are
# stored as text, matching the source data which uses dotted thousand
# separators rather than true numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '29.382.79'
$ws.Range("E2").Value = '  -1.71%  '

$ws.Range("D3").Value = '1.853.27'
$ws.Range("E3").Value = '  -1.15%  '

$ws.Range("D4").Value = '1.003'
$ws.Range("E4").Value = '  +0.38%  '

$ws.Range("D5").Value = '0.6999'
$ws.Range("E5").Value = '  -5.25%  '

$ws.Range("D6").Value = '238.91'
$ws.Range("E6").Value = '  -1.37%  '

$ws.Range("D7").Value = '1.002'
$ws.Range("E7").Value = '  +0.31%  '

$ws.Range("D8").Value = '0.3077'
$ws.Range("E8").Value = '  -2.57%  '

$ws.Range("D9").Value = '0.07389'
$ws.Range("E9").Value = '  +2.84%  '

$ws.Range("D10").Value = '23.72'
$ws.Range("E10").Value = '  -3.97%  '

$ws.Range("D11").Value = '0.08108'
$ws.Range("E11").Value = '  -3.17%  '

$ws.Range("D12").Value = '1.922.33'
$ws.Range("E12").Value = '  +2.09%  '

$ws.Range("D13").Value = '0.7278'
$ws.Range("E13").Value = '  -3.05%  '

$ws.Range("D14").Value = '5.220'
$ws.Range("E14").Value = '  -3.81%  '

$ws.Range("D15").Value = '89.92'
$ws.Range("E15").Value = '  -2.93%  '

$ws.Range("D16").Value = '29.557.16'
$ws.Range("E16").Value = '  -1.09%  '

$ws.Range("D17").Value = '5.928'
$ws.Range("E17").Value = '  -2.42%  '

$ws.Range("D18").Value = '243.12'
$ws.Range("E18").Value = '  -1.16%  '

$ws.Range("D19").Value = '0.000007764'
$ws.Range("E19").Value = '  -0.85%  '

$ws.Range("D20").Value = '13.18'
$ws.Range("E20").Value = '  -2.88%  '

$ws.Range("D21").Value = '1.002'
$ws.Range("E21").Value = '  +0.34%  '

$ws.Range("D22").Value = '2.133.76'
$ws.Range("E22").Value = '  +0.48%  '

$ws.Range("D23").Value = '1.003'
$ws.Range("E23").Value = '  +0.28%  '

$ws.Range("D24").Value = '7.636'
$ws.Range("E24").Value = '  -4.51%  '

$ws.Range("D25").Value = '0.1484'
$ws.Range("E25").Value = '  -4.18%  '

$ws.Range("B26").Value = 'Cosmos'
$ws.Range("C26").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D26").Value = '9.040'
$ws.Range("E26").Value = '  -2.39%  '

$ws.Range("B27").Value = 'Monero'
$ws.Range("C27").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D27").Value = '161.98'
$ws.Range("E27").Value = '  -1.80%  '

$ws.Range("D28").Value = '18.12'
$ws.Range("E28").Value = '  -2.82%  '

$ws.Range("D29").Value = '1.947'
$ws.Range("E29").Value = '  -4.29%  '

$ws.Range("D30").Value = '1.386'
$ws.Range("E30").Value = '  -8.10%  '

$ws.Range("D31").Value = '1.506'
$ws.Range("E31").Value = '  -1.84%  '

$ws.Range("D32").Value = '4.427'
$ws.Range("E32").Value = '  -3.60%  '

$ws.Range("D33").Value = '4.077'
$ws.Range("E33").Value = '  -4.86%  '

$ws.Range("D34").Value = '0.05309'
$ws.Range("E34").Value = '  -0.07%  '

$ws.Range("D35").Value = '1.201'
$ws.Range("E35").Value = '  -2.94%  '

$ws.Range("D36").Value = '0.7252'
$ws.Range("E36").Value = '  -3.90%  '

$ws.Range("D37").Value = '1.006'
$ws.Range("E37").Value = '  +0.71%  '

$ws.Range("D38").Value = '2.687'
$ws.Range("E38").Value = '  -0.03%  '

$ws.Range("D39").Value = '0.01869'
$ws.Range("E39").Value = '  -4.59%  '

$ws.Range("D40").Value = '2.704'
$ws.Range("E40").Value = '  -1.81%  '

$ws.Range("D41").Value = '0.4327'
$ws.Range("E41").Value = '  -3.92%  '

$ws.Range("D42").Value = '0.8734'
$ws.Range("E42").Value = '  +2.01%  '

$ws.Range("D43").Value = '5.918'
$ws.Range("E43").Value = '  -2.24%  '

$ws.Range("D44").Value = '70.35'
$ws.Range("E44").Value = '  -2.86%  '

$ws.Range("D45").Value = '1.002'
$ws.Range("E45").Value = '  +0.19%  '

$ws.Range("B46").Value = 'Maker'
$ws.Range("C46").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D46").Value = '1.027.58'
$ws.Range("E46").Value = '  -7.58%  '

$ws.Range("B47").Value = 'Quant'
$ws.Range("C47").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D47").Value = '102.57'
$ws.Range("E47").Value = '  -0.46%  '

$ws.Range("B48").Value = 'RocketPoolETH'
$ws.Range("C48").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D48").Value = '2.054.00'
$ws.Range("E48").Value = '  +1.58%  '

$ws.Range("B49").Value = 'Aptos'
$ws.Range("C49").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D49").Value = '7.284'
$ws.Range("E49").Value = '  -4.42%  '

$ws.Range("D50").Value = '1.747'
$ws.Range("E50").Value = '  -5.34%  '

$ws.Range("B51").Value = 'EnergySwap'
$ws.Range("C51").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D51").Value = '9.285'
$ws.Range("E51").Value = '  -2.16%  '
